$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 (pushes the existing "Ensemble" row down to row 7)
$ws.Rows.Item(6).Insert()

# Copy the formatting (bold font + border) from the label cell above onto the new A6 label cell
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$row2 = New-Object 'object[,]' 1,42
$row2[0,0] = 0.6616053383069038
$row2[0,1] = 0.03953544260148248
$row2[0,2] = 0.6557244632433089
$row2[0,3] = 0.6046864226291436
$row2[0,4] = 0.6636337423737134
$row2[0,5] = 0.6553153359261802
$row2[0,6] = 0.7286667273621734
$row2[0,7] = 0.6814570828139329
$row2[0,8] = 0.04100529591704264
$row2[0,9] = 0.6777850778042112
$row2[0,10] = 0.6796296017715363
$row2[0,11] = 0.6132348862837936
$row2[0,12] = 0.6957498402509451
$row2[0,13] = 0.7408860079591786
$row2[0,14] = 0.6904940944263334
$row2[0,15] = 0.02497588432533826
$row2[0,16] = 0.6700981445662297
$row2[0,17] = 0.6752796284534057
$row2[0,18] = 0.6834306249342563
$row2[0,19] = 0.6843439819720654
$row2[0,20] = 0.7393180922057093
$row2[0,21] = 0.8018820540695701
$row2[0,22] = 0.0352614642125757
$row2[0,23] = 0.8186480186480187
$row2[0,24] = 0.7318755459172916
$row2[0,25] = 0.8273710196771193
$row2[0,26] = 0.8155539649017638
$row2[0,27] = 0.8159617212036567
$row2[0,28] = 0.7873033124328401
$row2[0,29] = 0.03297090628898847
$row2[0,30] = 0.7890943372748754
$row2[0,31] = 0.727188940092166
$row2[0,32] = 0.8280731812989878
$row2[0,33] = 0.7971320171179919
$row2[0,34] = 0.7950280863801797
$row2[0,35] = 0.7719583017452321
$row2[0,36] = 0.04306347273586575
$row2[0,37] = 0.7738583246035472
$row2[0,38] = 0.6928501468582939
$row2[0,39] = 0.7721126899395151
$row2[0,40] = 0.8150644535581287
$row2[0,41] = 0.8059058937666748
$ws.Range("B2:AQ2").Value = $row2

$row3 = New-Object 'object[,]' 1,42
$row3[0,0] = 0.814105879425558
$row3[0,1] = 0.0453760117518966
$row3[0,2] = 0.8194782168186424
$row3[0,3] = 0.7761646390303554
$row3[0,4] = 0.8191047955140237
$row3[0,5] = 0.7630006916933911
$row3[0,6] = 0.8927810540713766
$row3[0,7] = 0.8165697171120989
$row3[0,8] = 0.03696497953419494
$row3[0,9] = 0.8003224472939523
$row3[0,10] = 0.7677276407115117
$row3[0,11] = 0.8165129939323488
$row3[0,12] = 0.8170451042371205
$row3[0,13] = 0.8812403993855606
$row3[0,14] = 0.7882568063073884
$row3[0,15] = 0.05587600468929411
$row3[0,16] = 0.7285415219345791
$row3[0,17] = 0.7577619238167203
$row3[0,18] = 0.8206778005165103
$row3[0,19] = 0.752367631741887
$row3[0,20] = 0.8819351535272449
$row3[0,21] = 0.86215461177632
$row3[0,22] = 0.03596446647499382
$row3[0,23] = 0.8705377220945462
$row3[0,24] = 0.8166666666666667
$row3[0,25] = 0.9147857424597663
$row3[0,26] = 0.8275846017781501
$row3[0,27] = 0.881198325882471
$row3[0,28] = 0.8669515214950227
$row3[0,29] = 0.03275762146175935
$row3[0,30] = 0.8608829416561643
$row3[0,31] = 0.8279620021555506
$row3[0,32] = 0.9149881539131611
$row3[0,33] = 0.8381095159320965
$row3[0,34] = 0.8928149938181409
$row3[0,35] = 0.8625959207273046
$row3[0,36] = 0.03022425430067465
$row3[0,37] = 0.8609885832900883
$row3[0,38] = 0.8279620021555506
$row3[0,39] = 0.9144298172415278
$row3[0,40] = 0.8381095159320965
$row3[0,41] = 0.8714896850172604
$ws.Range("B3:AQ3").Value = $row3

$row4 = New-Object 'object[,]' 1,42
$row4[0,0] = 0.8484254766264593
$row4[0,1] = 0.04609867125037707
$row4[0,2] = 0.840986818030491
$row4[0,3] = 0.7989640892866698
$row4[0,4] = 0.8715139511362958
$row4[0,5] = 0.8060902877541248
$row4[0,6] = 0.9245722369247152
$row4[0,7] = 0.8715984140059622
$row4[0,8] = 0.02798684160145205
$row4[0,9] = 0.8621038816486616
$row4[0,10] = 0.839348103864233
$row4[0,11] = 0.893185533104888
$row4[0,12] = 0.849143314620146
$row4[0,13] = 0.9142112367918819
$row4[0,14] = 0.8548517114911736
$row4[0,15] = 0.0336145026339547
$row4[0,16] = 0.8396574107139487
$row4[0,17] = 0.8087157764577119
$row4[0,18] = 0.8937036286740209
$row4[0,19] = 0.8384487490868707
$row4[0,20] = 0.8937329925233151
$row4[0,21] = 0.858065316010989
$row4[0,22] = 0.0337552764194182
$row4[0,23] = 0.828901872014974
$row4[0,24] = 0.8184368867480822
$row4[0,25] = 0.8702737617056016
$row4[0,26] = 0.8587319843138891
$row4[0,27] = 0.9139820752723978
$row4[0,28] = 0.8647697933549733
$row4[0,29] = 0.03275296086228412
$row4[0,30] = 0.8292586277147773
$row4[0,31] = 0.8281412136250846
$row4[0,32] = 0.8824142817651864
$row4[0,33] = 0.8700527683974202
$row4[0,34] = 0.9139820752723978
$row4[0,35] = 0.8625653863895216
$row4[0,36] = 0.03177525968687078
$row4[0,37] = 0.8295415195150263
$row4[0,38] = 0.8281412136250846
$row4[0,39] = 0.8711093551376792
$row4[0,40] = 0.8700527683974202
$row4[0,41] = 0.9139820752723978
$ws.Range("B4:AQ4").Value = $row4

$row5a = New-Object 'object[,]' 1,14
$row5a[0,0] = 0.8472833278025794
$row5a[0,1] = 0.03531711245032959
$row5a[0,2] = 0.8622178409412452
$row5a[0,3] = 0.7842929818476124
$row5a[0,4] = 0.8923855996043251
$row5a[0,5] = 0.8506432054819153
$row5a[0,6] = 0.8468770111377989
$row5a[0,7] = 0.8604651374710626
$row5a[0,8] = 0.03408157604204808
$row5a[0,9] = 0.8717925204812885
$row5a[0,10] = 0.7954964390448261
$row5a[0,11] = 0.8923855996043251
$row5a[0,12] = 0.8612331384866625
$row5a[0,13] = 0.8814179897382101
$ws.Range("B5:O5").Value = $row5a

$row5b = New-Object 'object[,]' 1,14
$row5b[0,0] = 0.825641043605575
$row5b[0,1] = 0.04659281028642469
$row5b[0,2] = 0.8710619275575197
$row5b[0,3] = 0.7598211574247751
$row5b[0,4] = 0.8833409069113157
$row5b[0,5] = 0.7922565641995499
$row5b[0,6] = 0.8217246619347144
$row5b[0,7] = 0.8367526725646955
$row5b[0,8] = 0.04220680691920378
$row5b[0,9] = 0.8710619275575197
$row5b[0,10] = 0.7715930859945006
$row5b[0,11] = 0.892902713922485
$row5b[0,12] = 0.8264809734142583
$row5b[0,13] = 0.8217246619347144
$ws.Range("W5:AJ5").Value = $row5b

$ws.Range("A6").Value = "RF"

$row6 = New-Object 'object[,]' 1,42
$row6[0,0] = 0.819964574942127
$row6[0,1] = 0.04345311618384834
$row6[0,2] = 0.7759252925486352
$row6[0,3] = 0.7701368405404772
$row6[0,4] = 0.8527884556110362
$row6[0,5] = 0.8182521109516017
$row6[0,6] = 0.8827201750588847
$row6[0,7] = 0.8240843164738699
$row6[0,8] = 0.02583954445945768
$row6[0,9] = 0.8138993814204805
$row6[0,10] = 0.8092031058578617
$row6[0,11] = 0.8286864918522833
$row6[0,12] = 0.7970157661026102
$row6[0,13] = 0.8716168371361133
$row6[0,14] = 0.7621559464350687
$row6[0,15] = 0.03101511909238578
$row6[0,16] = 0.7167428596040019
$row6[0,17] = 0.7434521605652489
$row6[0,18] = 0.8057484591709612
$row6[0,19] = 0.7604478265556364
$row6[0,20] = 0.7843884262794949
$row6[0,21] = 0.8047885220093512
$row6[0,22] = 0.031298577183661
$row6[0,23] = 0.756847819304923
$row6[0,24] = 0.7890694022909641
$row6[0,25] = 0.8084289703436425
$row6[0,26] = 0.8183941322580603
$row6[0,27] = 0.8512022858491665
$row6[0,28] = 0.8108420237142937
$row6[0,29] = 0.0362369166095132
$row6[0,30] = 0.7684288219977994
$row6[0,31] = 0.7994592458830108
$row6[0,32] = 0.8304785865269736
$row6[0,33] = 0.7850959463862689
$row6[0,34] = 0.8707475177774155
$row6[0,35] = 0.8191986420533903
$row6[0,36] = 0.0444742544506656
$row6[0,37] = 0.7717672373815133
$row6[0,38] = 0.7651082741065763
$row6[0,39] = 0.8384729239289437
$row6[0,40] = 0.8384081491042442
$row6[0,41] = 0.8822366257456736
$ws.Range("B6:AQ6").Value = $row6

$row7 = New-Object 'object[,]' 1,42
$row7[0,0] = 0.8564344881578118
$row7[0,1] = 0.04384619142444963
$row7[0,2] = 0.8305453571411018
$row7[0,3] = 0.8179406088324495
$row7[0,4] = 0.8922515656386625
$row7[0,5] = 0.8170073653944623
$row7[0,6] = 0.9244275437823826
$row7[0,7] = 0.865142642371238
$row7[0,8] = 0.04742723010437374
$row7[0,9] = 0.8286217619814977
$row7[0,10] = 0.8084052148568277
$row7[0,11] = 0.9038787082748531
$row7[0,12] = 0.8493515848354559
$row7[0,13] = 0.9354559419075548
$row7[0,14] = 0.8346262606246398
$row7[0,15] = 0.04031291308432402
$row7[0,16] = 0.8104153356030703
$row7[0,17] = 0.8205516108741916
$row7[0,18] = 0.8106010586910417
$row7[0,19] = 0.8166753185507876
$row7[0,20] = 0.9148879794041084
$row7[0,21] = 0.8817900478692774
$row7[0,22] = 0.05025721329436046
$row7[0,23] = 0.8923304374134681
$row7[0,24] = 0.7964198251674532
$row7[0,25] = 0.9356922371457135
$row7[0,26] = 0.8598090431292635
$row7[0,27] = 0.9246986964904892
$row7[0,28] = 0.8620746383400251
$row7[0,29] = 0.03439226246432351
$row7[0,30] = 0.8281377998102392
$row7[0,31] = 0.8164417326227941
$row7[0,32] = 0.9040359635648261
$row7[0,33] = 0.8705803506747646
$row7[0,34] = 0.8911773450275011
$row7[0,35] = 0.8627347404167841
$row7[0,36] = 0.04727350484673024
$row7[0,37] = 0.8930762530986496
$row7[0,38] = 0.7760878623357402
$row7[0,39] = 0.9146506683464736
$row7[0,40] = 0.8604164071906008
$row7[0,41] = 0.8694425111124562
$ws.Range("B7:AQ7").Value = $row7

